# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" column (E16:E52) listed the 37 mora periods from
# 1703 (Mar-2017) through 2003 (Mar-2020) in descending order. This
# update flips the listing to ascending order (1703 ... 2003) as part
# of refreshing the account-statement database.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$periods = @( `
    "1703","1704","1705","1706","1707","1708","1709","1710","1711","1712", `
    "1801","1802","1803","1804","1805","1806","1807","1808","1809","1810","1811","1812", `
    "1901","1902","1903","1904","1905","1906","1907","1908","1909","1910","1911","1912", `
    "2001","2002","2003" `
)

$firstRow = 16
for ($i = 0; $i -lt $periods.Length; $i++) {
    $ws.Cells.Item($firstRow + $i, 5).Value = $periods[$i]
}
